$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 5 (a "No Spaces Warning" QnA record), using the same
# column layout as rows 2-4: qid / question1 / question2 / Answer /
# markdown / ssml / topic / cardtitle / imageurl / displaytext1 /
# buttonvalue1 / displaytext2 / buttonvalue2.
$ws.Cells.Item(5,1).Value2  = "No Spaces Warning"
$ws.Cells.Item(5,2).Value2  = "What is Q and A Bot"
$ws.Cells.Item(5,3).Value2  = "What is QnaBot"
$ws.Cells.Item(5,4).Value2  = "The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Cells.Item(5,5).Value2  = "The Q and A Bot uses [Amazon Lex](https://aws.amazon.com/lex/) and [Alexa](https://developer.amazon.com/en-US/alexa) to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Cells.Item(5,6).Value2  = "<speak>The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer</speak>"
$ws.Cells.Item(5,7).Value2  = "Alexa"
$ws.Cells.Item(5,8).Value2  = "Alexa"
$ws.Cells.Item(5,10).Value2 = "Tell me about the Alexa Show."
$ws.Cells.Item(5,11).Value2 = "The Echo Show"
$ws.Cells.Item(5,12).Value2 = "Tell me about the Echo Dot"
$ws.Cells.Item(5,13).Value2 = "The Echo Dot"

# Re-apply the wrap-text formatting used by the other data rows for these
# columns.
$ws.Range("D5:F5").WrapText = $true
$ws.Range("H5").WrapText = $true
$ws.Range("J5:M5").WrapText = $true

# Column I (image url) is a real hyperlink, same target used by the other
# rows in the sheet.
$ws.Hyperlinks.Add($ws.Cells.Item(5,9), "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg")
$ws.Cells.Item(5,9).WrapText = $true

# Row height grows to fit the new wrapped content.
$ws.Rows.Item(5).RowHeight = 153

# Update the view so the new row is visible and selected, matching what
# the workbook looked like after the edit was made.
[void]$ws.Range("B5:M5").Select()

$wb.Save()
